$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: stash the two rows that need to move to a later position so we
#     don't clobber them while writing the newly uploaded word pairs. ---
$stashA2 = $ws.Cells.Item(2, 1).Value2   # "goodbye"
$stashB2 = $ws.Cells.Item(2, 2).Value2   # "adios"
$stashA3 = $ws.Cells.Item(3, 1).Value2   # "dog"
$stashB3 = $ws.Cells.Item(3, 2).Value2   # "perro"
$stashA4 = $ws.Cells.Item(4, 1).Value2   # "thank you"
$stashB4 = $ws.Cells.Item(4, 2).Value2   # "gracias"

$ws.Range("A2:B4").ClearContents()

# --- Step 2: write the newly uploaded word pairs straight into their final
#     row positions. The order below is the order the words were parsed
#     from the user's upload (only "=>" / "->" separators now allowed),
#     which is what determines shared-string table order. ---
$ws.Cells.Item(2, 1).Value = "one"
$ws.Cells.Item(2, 2).Value = "uno"

$ws.Cells.Item(4, 1).Value = "lifetime"
$ws.Cells.Item(4, 2).Value = "toda la vida"

$ws.Cells.Item(7, 1).Value = "see you later"
$ws.Cells.Item(7, 2).Value = "hasta luego"

$ws.Cells.Item(8, 1).Value = "how are you"
$ws.Cells.Item(8, 2).Value = "qué tal"

$ws.Cells.Item(5, 1).Value = "bank, bench"
$ws.Cells.Item(5, 2).Value = "banco"

$ws.Cells.Item(6, 1).Value = "light"
$ws.Cells.Item(6, 2).Value = "luz, ligero"

# --- Step 3: drop the stashed rows back in at their new positions. ---
$ws.Cells.Item(3, 1).Value = $stashA4
$ws.Cells.Item(3, 2).Value = $stashB4

$ws.Cells.Item(9, 1).Value = $stashA3
$ws.Cells.Item(9, 2).Value = $stashB3

$ws.Cells.Item(10, 1).Value = $stashA2
$ws.Cells.Item(10, 2).Value = $stashB2

# Row 11 is left blank (new trailing row), matching the final dimension A1:B11.

# --- Step 4: make formatting uniform (style index 1, "Segoe UI") across the
#     whole A1:B11 block, same as the rest of the table. Copy format only so
#     we don't disturb the values/shared-string table we just built. ---
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Range("A1:B11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 4b: every data row (including the trailing blank one) uses the
#     15pt row height the rest of the table uses. ---
$ws.Rows("1:11").RowHeight = 15

# --- Step 5: column widths for the two columns. ---
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 15

# --- Step 6: page setup / orientation. ---
$ws.PageSetup.Orientation = 1

# --- Step 7: restore the selection to match where the user's cursor ended
#     up (just past the last filled row). ---
$ws.Range("B16").Select() | Out-Null
